# "add user list to project"
#
# 1. The PI-hours sheet's cfop list for Romit Roy Choudhury's row (row 4) is
#    re-ordered (['cfop_CHOUDHURY', 'cfop_RRC'] -> ['cfop_RRC', 'cfop_CHOUDHURY']).
# 2. The project-hours sheet gains a new "users" column (E) listing the
#    user(s) associated with each project row.

$wb = $excel.ActiveWorkbook

# --- 1. PI hours: fix the cfop list ordering for row 4 (Romit Roy Choudhury) ---
$piSheet = $wb.Worksheets.Item("PI hours")
$piSheet.Range("G4").Value = "['cfop_RRC', 'cfop_CHOUDHURY']"

# --- 2. Project hours: add the "users" column ---
$projSheet = $wb.Worksheets.Item("project hours")

$projSheet.Range("E1").Value = "users"
# Match the bordered/bold/centered header style already used by B1:D1.
$projSheet.Range("D1").Copy()
$projSheet.Range("E1").PasteSpecial(-4122)

$users = @(
    "['Arun Lakshmanan']",
    "['Jonathan Hoff']",
    "['Ashutosh Dhekne']",
    "['Won Dong Shin']",
    "['Ashutosh Dhekne']",
    "['Harshal Maske']",
    "['Gavin Ananda']"
)

for ($i = 0; $i -lt $users.Length; $i++) {
    $row = $i + 2
    $projSheet.Cells.Item($row, 5).Value = $users[$i]
}
